$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: input renamed (nNodes -> nElements) and new measurement values ---
$ws1.Range("A9").Value = "nElements"
$ws1.Range("B8").Value = 0.008
$ws1.Range("B9").Value = 10

# --- Sheet1: recompute dTau as the average of v0/v1 instead of the previous extrapolation ---
$ws1.Range("F2").Formula = "= (B6 + B7)/2"
$ws1.Range("F3").Formula = "=B5/F2"

# --- Sheet1 cosmetics: widen column F, keep it custom (not best-fit) ---
$ws1.Range("F1:F9").EntireColumn.ColumnWidth = 15

# --- Duplicate Sheet1 twice so each measurement series gets its own tab for the thesis graphs ---
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item($ws1.Index + 1)
$ws1.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item($ws2.Index + 1)

# --- Sheet1 (2): second measurement series ---
$ws2.Range("B6").Value = 0.005
$ws2.Range("B7").Value = 0.0055
$ws2.Range("B8").Value = 0.0025
$ws2.Range("B9").Value = 15
$ws2.Range("B23").Select()

# --- Sheet1 (3): third measurement series ---
$ws3.Range("B6").Value = 0.005
$ws3.Range("B7").Value = 0.0055
$ws3.Range("B8").Value = 0.005
$ws3.Range("B9").Value = 5
$ws3.Range("B8").Select()

# --- Restore Sheet1 as the active tab/selection ---
$ws1.Activate()
$ws1.Range("B9").Select()
